$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
# Row 28
$ws.Range("H28").Value = 1676.4615
$ws.Range("I28").Value = 1791.1666
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 1791.1666
$ws.Range("L28").Value = 300
$ws.Range("M28").Value = -1306.1666
$ws.Range("N28").Value = -1270
# Row 40
$ws.Range("H40").Value = 2293.1428
$ws.Range("I40").Value = 2190.4
$ws.Range("J40").Value = 2550
$ws.Range("K40").Value = 2190.4
$ws.Range("L40").Value = 2550
$ws.Range("M40").Value = -2015.4
$ws.Range("N40").Value = -2900
# Row 62
$ws.Range("H62").Value = 37041036
$ws.Range("I62").Value = 55560556
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 55560556
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -55559932
$ws.Range("N62").Value = -3248
# Row 65
$ws.Range("H65").Value = 37041036
$ws.Range("I65").Value = 55560556
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 277802780
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -277799660
$ws.Range("N65").Value = -16240
# Row 98
$ws.Range("H98").Value = 3701.8708
$ws.Range("I98").Value = 3952.0386
$ws.Range("J98").Value = 2401
$ws.Range("K98").Value = 3952.0386
$ws.Range("L98").Value = 2401
$ws.Range("M98").Value = -2454.0386
$ws.Range("N98").Value = -5397
# Row 113
$ws.Range("H113").Value = 2006
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
# Row 118
$ws.Range("H118").Value = 951.5
$ws.Range("I118").Value = 345.42856
$ws.Range("K118").Value = 1036.28568
$ws.Range("M118").Value = 620.71432
# Row 122
$ws.Range("H122").Value = 3701.8708
$ws.Range("I122").Value = 3952.0386
$ws.Range("J122").Value = 2401
$ws.Range("K122").Value = 11856.1158
$ws.Range("L122").Value = 7203
$ws.Range("M122").Value = -9406.1158
$ws.Range("N122").Value = -12103
# Row 135
$ws.Range("H135").Value = 29412340
$ws.Range("I135").Value = 372.46155
$ws.Range("J135").Value = 125001230
$ws.Range("K135").Value = 3352.15395
$ws.Range("L135").Value = 1125011070
$ws.Range("M135").Value = -817.1539499999999
$ws.Range("N135").Value = -1125016140
# Row 138
$ws.Range("H138").Value = 589200.8
$ws.Range("I138").Value = 1480.8125
$ws.Range("J138").Value = 751330.4399999999
$ws.Range("K138").Value = 4442.4375
$ws.Range("L138").Value = 2253991.32
$ws.Range("M138").Value = 697.5625
$ws.Range("N138").Value = -2264271.32

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4857.763
$ws.Range("I32").Value = 4857.763
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4857.763
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4570.763
$ws.Range("N32").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1334.9286
$ws.Range("I107").Value = 1093.75
$ws.Range("J107").Value = 1656.5
$ws.Range("K107").Value = 1093.75
$ws.Range("L107").Value = 1656.5
$ws.Range("M107").Value = 826.25
$ws.Range("N107").Value = -5496.5
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1466.3
$ws.Range("I31").Value = 1642.2222
$ws.Range("J31").Value = 1390.9048
$ws.Range("K31").Value = 1642.2222
$ws.Range("L31").Value = 1390.9048
$ws.Range("M31").Value = -1347.2222
$ws.Range("N31").Value = -1980.9048
# Row 34
$ws.Range("H34").Value = 1466.3
$ws.Range("I34").Value = 1642.2222
$ws.Range("J34").Value = 1390.9048
$ws.Range("K34").Value = 1642.2222
$ws.Range("L34").Value = 1390.9048
$ws.Range("M34").Value = -1440.2222
$ws.Range("N34").Value = -1794.9048
# Row 58
$ws.Range("H58").Value = 1696.591
$ws.Range("I58").Value = 1401.6428
$ws.Range("J58").Value = 2212.75
$ws.Range("K58").Value = 1401.6428
$ws.Range("L58").Value = 2212.75
$ws.Range("M58").Value = -1198.6428
$ws.Range("N58").Value = -2618.75
# Row 95
$ws.Range("H95").Value = 10002.4
$ws.Range("J95").Value = 10002.4
$ws.Range("L95").Value = 10002.4
$ws.Range("N95").Value = -15494.4
# Row 132
$ws.Range("H132").Value = 1934.0952
$ws.Range("I132").Value = 1404.5385
$ws.Range("J132").Value = 2794.625
$ws.Range("K132").Value = 4213.6155
$ws.Range("L132").Value = 8383.875
$ws.Range("M132").Value = -1683.6155
$ws.Range("N132").Value = -13443.875
# Row 136
$ws.Range("H136").Value = 1696.591
$ws.Range("I136").Value = 1401.6428
$ws.Range("J136").Value = 2212.75
$ws.Range("K136").Value = 4204.928400000001
$ws.Range("L136").Value = 6638.25
$ws.Range("M136").Value = -1654.928400000001
$ws.Range("N136").Value = -11738.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 224.26315
$ws.Range("I14").Value = 224.26315
$ws.Range("K14").Value = 672.78945
$ws.Range("M14").Value = -499.78945
# Row 32
$ws.Range("H32").Value = 1717
$ws.Range("J32").Value = 2300
$ws.Range("L32").Value = 6900
$ws.Range("N32").Value = -7466
# Row 39
$ws.Range("H39").Value = 2945
$ws.Range("J39").Value = 2855
$ws.Range("L39").Value = 8565
$ws.Range("N39").Value = -9153
# Row 113
$ws.Range("H113").Value = 658.9143
$ws.Range("I113").Value = 584.5
$ws.Range("J113").Value = 680.96295
$ws.Range("K113").Value = 1753.5
$ws.Range("L113").Value = 2042.88885
$ws.Range("M113").Value = 416.5
$ws.Range("N113").Value = -6382.888849999999
# Row 122
$ws.Range("H122").Value = 1103.8889
$ws.Range("J122").Value = 1136.875
$ws.Range("L122").Value = 10231.875
$ws.Range("N122").Value = -15131.875
# Row 132
$ws.Range("H132").Value = 791.1667
$ws.Range("I132").Value = 791.1667
$ws.Range("K132").Value = 7120.5003
$ws.Range("M132").Value = -4590.5003
# Row 140
$ws.Range("H140").Value = 23263.469
$ws.Range("I140").Value = 52006.35
$ws.Range("K140").Value = 156019.05
$ws.Range("M140").Value = -150839.05

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 8066325
$ws.Range("I122").Value = 1868.1904
$ws.Range("J122").Value = 25001684
$ws.Range("K122").Value = 5604.5712
$ws.Range("L122").Value = 75005052
$ws.Range("M122").Value = -3154.5712
$ws.Range("N122").Value = -75009952

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2255.92
$ws.Range("I40").Value = 2094.9048
$ws.Range("J40").Value = 3101.25
$ws.Range("K40").Value = 2094.9048
$ws.Range("L40").Value = 3101.25
$ws.Range("M40").Value = -1958.9048
$ws.Range("N40").Value = -3373.25
# Row 46
$ws.Range("H46").Value = 5362.222
# Row 93
$ws.Range("H93").Value = 1000.1429
$ws.Range("I93").Value = 999.4
$ws.Range("J93").Value = 1002
$ws.Range("K93").Value = 999.4
$ws.Range("L93").Value = 1002
$ws.Range("M93").Value = 248.6
$ws.Range("N93").Value = -3498

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 10000
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("N69").Value = -11498
# Row 72
$ws.Range("H72").Value = 10000
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 30000
$ws.Range("N72").Value = -37488
# Row 98
$ws.Range("H98").Value = 21750
$ws.Range("J98").Value = 21750
$ws.Range("L98").Value = 21750
$ws.Range("N98").Value = -27740
# Row 107
$ws.Range("H107").Value = 569.2857
$ws.Range("I107").Value = 472
$ws.Range("J107").Value = 812.5
$ws.Range("K107").Value = 1416
$ws.Range("L107").Value = 2437.5
$ws.Range("M107").Value = 504
$ws.Range("N107").Value = -6277.5
# Row 136
$ws.Range("H136").Value = 906.60974
$ws.Range("I136").Value = 834.7742
$ws.Range("J136").Value = 1129.3
$ws.Range("K136").Value = 2504.3226
$ws.Range("L136").Value = 3387.9
$ws.Range("M136").Value = 45.67740000000003
$ws.Range("N136").Value = -8487.9

